# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin name / link swapped between rows 39 and 40 (re-sorted ranking).
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'

# Updated Price (column D) and Volume(1h) (column E) values per row.
$priceVolumeUpdates = @(
    @{ Cell = 'D2'; Value = '24.248.36' },
    @{ Cell = 'E2'; Value = '  +14.24%  ' },
    @{ Cell = 'D3'; Value = '1.674.69' },
    @{ Cell = 'D4'; Value = '0.9993' },
    @{ Cell = 'E4'; Value = '  -0.02%  ' },
    @{ Cell = 'D5'; Value = '308.70' },
    @{ Cell = 'E5'; Value = '  +9.71%  ' },
    @{ Cell = 'D6'; Value = '0.9964' },
    @{ Cell = 'E6'; Value = '  +3.89%  ' },
    @{ Cell = 'D7'; Value = '0.3735' },
    @{ Cell = 'E7'; Value = '  +3.27%  ' },
    @{ Cell = 'D8'; Value = '0.3439' },
    @{ Cell = 'E8'; Value = '  +8.35%  ' },
    @{ Cell = 'D9'; Value = '48.01' },
    @{ Cell = 'E9'; Value = '  +17.74%  ' },
    @{ Cell = 'D10'; Value = '1.188' },
    @{ Cell = 'E10'; Value = '  +7.12%  ' },
    @{ Cell = 'D11'; Value = '0.07298' },
    @{ Cell = 'E11'; Value = '  +7.62%  ' },
    @{ Cell = 'D12'; Value = '1.001' },
    @{ Cell = 'E12'; Value = '  +0.68%  ' },
    @{ Cell = 'D13'; Value = '20.58' },
    @{ Cell = 'E13'; Value = '  +9.91%  ' },
    @{ Cell = 'D14'; Value = '6.076' },
    @{ Cell = 'E14'; Value = '  +7.80%  ' },
    @{ Cell = 'D15'; Value = '6.770' },
    @{ Cell = 'E15'; Value = '  +6.88%  ' },
    @{ Cell = 'D16'; Value = '1.670.35' },
    @{ Cell = 'E16'; Value = '  +9.07%  ' },
    @{ Cell = 'E17'; Value = '  +6.52%  ' },
    @{ Cell = 'D18'; Value = '0.9975' },
    @{ Cell = 'E18'; Value = '  +3.90%  ' },
    @{ Cell = 'D19'; Value = '0.06722' },
    @{ Cell = 'E19'; Value = '  +11.35%  ' },
    @{ Cell = 'D20'; Value = '82.07' },
    @{ Cell = 'E20'; Value = '  +14.47%  ' },
    @{ Cell = 'D21'; Value = '16.49' },
    @{ Cell = 'E21'; Value = '  +9.97%  ' },
    @{ Cell = 'E22'; Value = '  +9.09%  ' },
    @{ Cell = 'E23'; Value = '  +6.28%  ' },
    @{ Cell = 'D24'; Value = '24.150.16' },
    @{ Cell = 'E24'; Value = '  +13.82%  ' },
    @{ Cell = 'D25'; Value = '2.406' },
    @{ Cell = 'E25'; Value = '  +4.46%  ' },
    @{ Cell = 'D26'; Value = '3.385' },
    @{ Cell = 'E26'; Value = '  -8.03%  ' },
    @{ Cell = 'D27'; Value = '2.671' },
    @{ Cell = 'E27'; Value = '  +21.77%  ' },
    @{ Cell = 'D28'; Value = '152.01' },
    @{ Cell = 'E28'; Value = '  +3.13%  ' },
    @{ Cell = 'D29'; Value = '19.52' },
    @{ Cell = 'E29'; Value = '  +10.51%  ' },
    @{ Cell = 'D30'; Value = '1.852.82' },
    @{ Cell = 'E30'; Value = '  +8.95%  ' },
    @{ Cell = 'D31'; Value = '127.61' },
    @{ Cell = 'E31'; Value = '  +8.68%  ' },
    @{ Cell = 'D32'; Value = '6.324' },
    @{ Cell = 'E32'; Value = '  +22.72%  ' },
    @{ Cell = 'D33'; Value = '4.028' },
    @{ Cell = 'E33'; Value = '  -2.25%  ' },
    @{ Cell = 'D34'; Value = '0.9876' },
    @{ Cell = 'E34'; Value = '  +16.38%  ' },
    @{ Cell = 'D35'; Value = '1.746' },
    @{ Cell = 'E35'; Value = '  +16.38%  ' },
    @{ Cell = 'D36'; Value = '0.08430' },
    @{ Cell = 'E36'; Value = '  +5.25%  ' },
    @{ Cell = 'D37'; Value = '12.38' },
    @{ Cell = 'E37'; Value = '  +15.53%  ' },
    @{ Cell = 'D38'; Value = '8.948' },
    @{ Cell = 'E38'; Value = '  +16.91%  ' },
    @{ Cell = 'D39'; Value = '5.367' },
    @{ Cell = 'E39'; Value = '  +9.44%  ' },
    @{ Cell = 'D40'; Value = '0.06408' },
    @{ Cell = 'E40'; Value = '  +9.59%  ' },
    @{ Cell = 'D41'; Value = '1.296' },
    @{ Cell = 'E41'; Value = '  +6.23%  ' },
    @{ Cell = 'D42'; Value = '0.02355' },
    @{ Cell = 'E42'; Value = '  +12.76%  ' },
    @{ Cell = 'D43'; Value = '0.2119' },
    @{ Cell = 'E43'; Value = '  +11.20%  ' },
    @{ Cell = 'D44'; Value = '0.6143' },
    @{ Cell = 'E44'; Value = '  +13.19%  ' },
    @{ Cell = 'D45'; Value = '0.9960' },
    @{ Cell = 'E45'; Value = '  +3.79%  ' },
    @{ Cell = 'D46'; Value = '3.807' },
    @{ Cell = 'E46'; Value = '  +7.20%  ' },
    @{ Cell = 'D47'; Value = '13.21' },
    @{ Cell = 'E47'; Value = '  +6.05%  ' },
    @{ Cell = 'D48'; Value = '0.5984' },
    @{ Cell = 'E48'; Value = '  +10.35%  ' },
    @{ Cell = 'D49'; Value = '127.07' },
    @{ Cell = 'E49'; Value = '  +4.42%  ' },
    @{ Cell = 'D50'; Value = '2.028' },
    @{ Cell = 'E50'; Value = '  +8.77%  ' },
    @{ Cell = 'D51'; Value = '0.07128' },
    @{ Cell = 'E51'; Value = '  +7.87%  ' }
)

foreach ($update in $priceVolumeUpdates) {
    $range = $ws.Range($update.Cell)
    if ($update.Cell.StartsWith("D")) {
        # Force text format so numeric-looking prices (e.g. '0.9993')
        # are stored as text, matching the source data, not coerced to numbers.
        $range.NumberFormat = '@'
    }
    $range.Value = $update.Value
}
